$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$tr.Paragraphs(6, 1).Text = "Four key components: physical model, virtual model, data, and communication services."
$tr.Paragraphs(7, 1).Text = "Currently mostly standalone systems with limited interoperability, especially in collected and exploited data."
